# Controlo de progresso semana 10
# Update the weekly progress-tracking report: new report date, and updated
# progress percentages for a few tasks in the "Avaliação"/"Instalação"
# sections of the Gantt-style tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report date (A2): moved one week forward ---
$ws.Range("A2").Value = 45665

# --- Row 36 (T5.1 - Relatorio) ---
# Previously-reported progress (F36) and current progress (G36) both
# increase; F36 also picks up G36's cell formatting (the thin/medium
# border combination used by the rest of column G), matching how the
# task's "already-completed" portion now shares the same look.
$ws.Range("G36").Copy()
$ws.Range("F36").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F36").Value = 0.3
$ws.Range("G36").Value = 0.5

# --- Rows 39-41 (T5.4, T5.5, T5.6) marked fully complete ---
# Previously-reported progress (F) catches up to the already-complete
# current progress (G), and F adopts G's border formatting for each row.
foreach ($row in 39..41) {
    $src = $ws.Range("G$row")
    $dst = $ws.Range("F$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)           # xlPasteFormats
    $dst.Value = 1
}

$excel.CutCopyMode = $false

# --- Leave the cursor where the author ended up editing ---
[void]$ws.Range("I36").Select()
